$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the SCENARIO_DESC (column D) text for rows 2-5, using the new
# short-form labels instead of the previous step-by-step instructions.
$ws.Range("D2").Value = "Tambah Setup Emiten"
$ws.Range("D3").Value = "View Setup Emiten"
$ws.Range("D4").Value = "Ubah Setup Emiten"
$ws.Range("D5").Value = "Hapus Setup Emiten"

# Adjust row heights to match the new, shorter text content.
$ws.Rows.Item(2).RowHeight = 30
$ws.Rows.Item(3).RowHeight = 30
$ws.Rows.Item(4).AutoFit()
$ws.Rows.Item(5).RowHeight = 30

# Update the selected cell/range in the sheet view.
$ws.Range("D5").Select()
